$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 257.42856
$ws.Range("I33").Value = 216.66667
$ws.Range("J33").Value = 502
$ws.Range("K33").Value = 216.66667
$ws.Range("L33").Value = 502
$ws.Range("M33").Value = 12.33332999999999
$ws.Range("N33").Value = -960
$ws.Range("H69").Value = 40000.91
$ws.Range("I69").Value = 118505
$ws.Range("J69").Value = 22555.555
$ws.Range("K69").Value = 355515
$ws.Range("L69").Value = 67666.66500000001
$ws.Range("M69").Value = -354641
$ws.Range("N69").Value = -69414.66500000001
$ws.Range("H72").Value = 40000.91
$ws.Range("I72").Value = 118505
$ws.Range("J72").Value = 22555.555
$ws.Range("K72").Value = 1066545
$ws.Range("L72").Value = 202999.995
$ws.Range("M72").Value = -1062177
$ws.Range("N72").Value = -211735.995
$ws.Range("I125").Value = 250001500
$ws.Range("J125").Value = 100018900
$ws.Range("K125").Value = 2250013500
$ws.Range("L125").Value = 900170100
$ws.Range("M125").Value = -2250011040
$ws.Range("N125").Value = -900175020
$ws.Range("H137").Value = 1575.2858
$ws.Range("I137").Value = 1088.2222
$ws.Range("J137").Value = 1940.5834
$ws.Range("K137").Value = 3264.6666
$ws.Range("L137").Value = 5821.7502
$ws.Range("M137").Value = -714.6665999999996
$ws.Range("N137").Value = -10921.7502
$ws.Range("H138").Value = 4157.028
$ws.Range("I138").Value = 2039.6428
$ws.Range("J138").Value = 5504.4546
$ws.Range("K138").Value = 6118.928400000001
$ws.Range("L138").Value = 16513.3638
$ws.Range("M138").Value = -978.9284000000007
$ws.Range("N138").Value = -26793.3638
$ws.Range("H141").Value = 2283.4783
$ws.Range("I141").Value = 2220.1667
$ws.Range("J141").Value = 2511.4
$ws.Range("K141").Value = 6660.500100000001
$ws.Range("L141").Value = 7534.200000000001
$ws.Range("M141").Value = -1480.500100000001
$ws.Range("N141").Value = -17894.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 253070
$ws.Range("I6").Value = 253070
$ws.Range("K6").Value = 253070
$ws.Range("M6").Value = -252897
$ws.Range("H10").Value = 669316.7
$ws.Range("I10").Value = 669316.7
$ws.Range("K10").Value = 669316.7
$ws.Range("M10").Value = -669146.7
$ws.Range("H122").Value = 1430.1818
$ws.Range("I122").Value = 1223.2
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 3669.6
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -1219.6
$ws.Range("N122").Value = -15400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 857
$ws.Range("J64").Value = 789.75
$ws.Range("L64").Value = 789.75
$ws.Range("N64").Value = -1239.75
$ws.Range("H67").Value = 857
$ws.Range("J67").Value = 789.75
$ws.Range("L67").Value = 789.75
$ws.Range("N67").Value = -2349.75
$ws.Range("H105").Value = 2710
$ws.Range("I105").Value = 2833.3333
$ws.Range("J105").Value = 2525
$ws.Range("K105").Value = 2833.3333
$ws.Range("L105").Value = 2525
$ws.Range("M105").Value = -1086.3333
$ws.Range("N105").Value = -6019

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 373.25
$ws.Range("I22").Value = 446.5
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 446.5
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -96.5
$ws.Range("N22").Value = -1000
$ws.Range("H86").Value = 16248.75
$ws.Range("I86").Value = 17398.2
$ws.Range("J86").Value = 14333
$ws.Range("K86").Value = 17398.2
$ws.Range("L86").Value = 14333
$ws.Range("M86").Value = -16275.2
$ws.Range("N86").Value = -16579
$ws.Range("H89").Value = 16248.75
$ws.Range("I89").Value = 17398.2
$ws.Range("J89").Value = 14333
$ws.Range("K89").Value = 86991
$ws.Range("L89").Value = 71665
$ws.Range("M89").Value = -81375
$ws.Range("N89").Value = -82897
$ws.Range("H94").Value = 88482
$ws.Range("I94").Value = 160822.58
$ws.Range("J94").Value = 4084.6667
$ws.Range("K94").Value = 160822.58
$ws.Range("L94").Value = 4084.6667
$ws.Range("M94").Value = -160371.58
$ws.Range("N94").Value = -4986.6667
$ws.Range("H134").Value = 3127.111
$ws.Range("I134").Value = 3028.5715
$ws.Range("K134").Value = 9085.7145
$ws.Range("M134").Value = -6550.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 333
$ws.Range("I5").Value = 333
$ws.Range("K5").Value = 999
$ws.Range("M5").Value = -887
$ws.Range("H113").Value = 1754.3334
$ws.Range("I113").Value = 1446.3334
$ws.Range("K113").Value = 4339.0002
$ws.Range("M113").Value = -2169.0002
$ws.Range("H132").Value = 3693.375
$ws.Range("J132").Value = 3258
$ws.Range("L132").Value = 29322
$ws.Range("N132").Value = -34382
$ws.Range("H135").Value = 333
$ws.Range("I135").Value = 333
$ws.Range("K135").Value = 2997
$ws.Range("M135").Value = -462
$ws.Range("H137").Value = 2432.1667
$ws.Range("I137").Value = 1398.5
$ws.Range("J137").Value = 4499.5
$ws.Range("K137").Value = 4195.5
$ws.Range("L137").Value = 13498.5
$ws.Range("M137").Value = 904.5
$ws.Range("N137").Value = -23698.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4998.5
$ws.Range("I68").Value = 4998.5
$ws.Range("K68").Value = 4998.5
$ws.Range("M68").Value = -4249.5
$ws.Range("H71").Value = 4998.5
$ws.Range("I71").Value = 4998.5
$ws.Range("K71").Value = 24992.5
$ws.Range("M71").Value = -21248.5
$ws.Range("H106").Value = 12752.714
$ws.Range("J106").Value = 12752.714
$ws.Range("L106").Value = 12752.714
$ws.Range("N106").Value = -15276.714
$ws.Range("H136").Value = 2738.8333
$ws.Range("I136").Value = 1886.8
$ws.Range("J136").Value = 6999
$ws.Range("K136").Value = 5660.4
$ws.Range("L136").Value = 20997
$ws.Range("M136").Value = -3110.4
$ws.Range("N136").Value = -26097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19840.8
$ws.Range("I62").Value = 27501.334
$ws.Range("J62").Value = 8350
$ws.Range("K62").Value = 27501.334
$ws.Range("L62").Value = 8350
$ws.Range("M62").Value = -26877.334
$ws.Range("N62").Value = -9598
$ws.Range("H65").Value = 19840.8
$ws.Range("I65").Value = 27501.334
$ws.Range("J65").Value = 8350
$ws.Range("K65").Value = 137506.67
$ws.Range("L65").Value = 41750
$ws.Range("M65").Value = -134386.67
$ws.Range("N65").Value = -47990
$ws.Range("H126").Value = 4236.5405
$ws.Range("I126").Value = 4079.913
$ws.Range("J126").Value = 4493.857
$ws.Range("K126").Value = 12239.739
$ws.Range("L126").Value = 13481.571
$ws.Range("M126").Value = -9769.739
$ws.Range("N126").Value = -18421.571
$ws.Range("H132").Value = 7000.5
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060
